# Auto-generated edit script: update market-price derived columns (H-N)
# on the Aegis_Profits leve-profit tables, one block per worksheet tab.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value2 = 619.9
$ws.Range("I2").Value2 = 628.4286
$ws.Range("J2").Value2 = 600
$ws.Range("K2").Value2 = 628.4286
$ws.Range("L2").Value2 = 600
$ws.Range("M2").Value2 = -515.4286
$ws.Range("N2").Value2 = -826
$ws.Range("H18").Value2 = 7132.4
$ws.Range("J18").Value2 = 50251
$ws.Range("L18").Value2 = 50251
$ws.Range("N18").Value2 = -50819
$ws.Range("H80").Value2 = 77441.30499999999
$ws.Range("I80").Value2 = 111446.336
$ws.Range("J80").Value2 = 59438.65
$ws.Range("K80").Value2 = 334339.008
$ws.Range("L80").Value2 = 178315.95
$ws.Range("M80").Value2 = -333341.008
$ws.Range("N80").Value2 = -180311.95
$ws.Range("H83").Value2 = 77441.30499999999
$ws.Range("I83").Value2 = 111446.336
$ws.Range("J83").Value2 = 59438.65
$ws.Range("K83").Value2 = 1003017.024
$ws.Range("L83").Value2 = 534947.85
$ws.Range("M83").Value2 = -998025.024
$ws.Range("N83").Value2 = -544931.85
$ws.Range("H88").Value2 = 2628.5715
$ws.Range("I88").Value2 = 1060.6
$ws.Range("J88").Value2 = 3118.5625
$ws.Range("K88").Value2 = 1060.6
$ws.Range("L88").Value2 = 3118.5625
$ws.Range("M88").Value2 = -654.5999999999999
$ws.Range("N88").Value2 = -3930.5625
$ws.Range("H91").Value2 = 2628.5715
$ws.Range("I91").Value2 = 1060.6
$ws.Range("J91").Value2 = 3118.5625
$ws.Range("K91").Value2 = 1060.6
$ws.Range("L91").Value2 = 3118.5625
$ws.Range("M91").Value2 = 343.4000000000001
$ws.Range("N91").Value2 = -5926.5625
$ws.Range("H129").Value2 = 483293.25
$ws.Range("J129").Value2 = 515352
$ws.Range("L129").Value2 = 1546056
$ws.Range("N129").Value2 = -1556056
$ws.Range("H138").Value2 = 2624
$ws.Range("I138").Value2 = 2611.4167
$ws.Range("J138").Value2 = 2626.4753
$ws.Range("K138").Value2 = 7834.250100000001
$ws.Range("L138").Value2 = 7879.4259
$ws.Range("M138").Value2 = -2694.250100000001
$ws.Range("N138").Value2 = -18159.4259

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value2 = 0
$ws.Range("J51").Value2 = 0
$ws.Range("L51").Value2 = 0
$ws.Range("H88").Value2 = 3385.9
$ws.Range("J88").Value2 = 2575.8333
$ws.Range("L88").Value2 = 2575.8333
$ws.Range("N88").Value2 = -3387.8333
$ws.Range("H91").Value2 = 3385.9
$ws.Range("J91").Value2 = 2575.8333
$ws.Range("L91").Value2 = 2575.8333
$ws.Range("N91").Value2 = -5383.8333
$ws.Range("H122").Value2 = 4980
$ws.Range("I122").Value2 = 4980
$ws.Range("J122").Value2 = 0
$ws.Range("K122").Value2 = 14940
$ws.Range("L122").Value2 = 0
$ws.Range("M122").Value2 = -12490
$ws.Range("H132").Value2 = 1997.5555
$ws.Range("I132").Value2 = 1420.3684
$ws.Range("K132").Value2 = 4261.1052
$ws.Range("M132").Value2 = -1731.1052
$ws.Range("N51").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 40890.71
$ws.Range("I86").Value2 = 51936.168
$ws.Range("K86").Value2 = 51936.168
$ws.Range("M86").Value2 = -50813.168
$ws.Range("H89").Value2 = 40890.71
$ws.Range("I89").Value2 = 51936.168
$ws.Range("K89").Value2 = 259680.84
$ws.Range("M89").Value2 = -254064.84

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value2 = 10400
$ws.Range("I14").Value2 = 1000
$ws.Range("K14").Value2 = 1000
$ws.Range("M14").Value2 = -830
$ws.Range("H41").Value2 = 9267.143
$ws.Range("I41").Value2 = 3693.75
$ws.Range("J41").Value2 = 12696.923
$ws.Range("K41").Value2 = 3693.75
$ws.Range("L41").Value2 = 12696.923
$ws.Range("M41").Value2 = -3265.75
$ws.Range("N41").Value2 = -13552.923
$ws.Range("H50").Value2 = 14568
$ws.Range("J50").Value2 = 14568
$ws.Range("L50").Value2 = 14568
$ws.Range("N50").Value2 = -15818
$ws.Range("H51").Value2 = 7952.7144
$ws.Range("J51").Value2 = 7929.8335
$ws.Range("L51").Value2 = 7929.8335
$ws.Range("N51").Value2 = -9401.833500000001
$ws.Range("H59").Value2 = 21836.666
$ws.Range("J59").Value2 = 21836.666
$ws.Range("L59").Value2 = 21836.666
$ws.Range("N59").Value2 = -24126.666
$ws.Range("H60").Value2 = 10957.392
$ws.Range("I60").Value2 = 8000
$ws.Range("J60").Value2 = 11091.818
$ws.Range("K60").Value2 = 8000
$ws.Range("L60").Value2 = 11091.818
$ws.Range("M60").Value2 = -7489
$ws.Range("N60").Value2 = -12113.818
$ws.Range("H61").Value2 = 7952.7144
$ws.Range("J61").Value2 = 7929.8335
$ws.Range("L61").Value2 = 7929.8335
$ws.Range("N61").Value2 = -8625.833500000001
$ws.Range("H74").Value2 = 22974.572
$ws.Range("J74").Value2 = 22974.572
$ws.Range("L74").Value2 = 22974.572
$ws.Range("N74").Value2 = -24722.572
$ws.Range("H77").Value2 = 22974.572
$ws.Range("J77").Value2 = 22974.572
$ws.Range("L77").Value2 = 68923.716
$ws.Range("N77").Value2 = -77659.716
$ws.Range("H127").Value2 = 0
$ws.Range("J127").Value2 = 0
$ws.Range("L127").Value2 = 0
$ws.Range("H132").Value2 = 4364.647
$ws.Range("I132").Value2 = 4454.087
$ws.Range("J132").Value2 = 4177.636
$ws.Range("K132").Value2 = 13362.261
$ws.Range("L132").Value2 = 12532.908
$ws.Range("M132").Value2 = -10832.261
$ws.Range("N132").Value2 = -17592.908
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value2 = 776.41
$ws.Range("I131").Value2 = 398.6111
$ws.Range("J131").Value2 = 859.3415
$ws.Range("K131").Value2 = 1195.8333
$ws.Range("L131").Value2 = 2578.0245
$ws.Range("M131").Value2 = 3844.1667
$ws.Range("N131").Value2 = -12658.0245

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 100002360
$ws.Range("I80").Value2 = 500002500
$ws.Range("J80").Value2 = 2329.25
$ws.Range("K80").Value2 = 500002500
$ws.Range("L80").Value2 = 2329.25
$ws.Range("M80").Value2 = -500001502
$ws.Range("N80").Value2 = -4325.25
$ws.Range("H83").Value2 = 100002360
$ws.Range("I83").Value2 = 500002500
$ws.Range("J83").Value2 = 2329.25
$ws.Range("K83").Value2 = 2500012500
$ws.Range("L83").Value2 = 11646.25
$ws.Range("M83").Value2 = -2500007508
$ws.Range("N83").Value2 = -21630.25
$ws.Range("H122").Value2 = 0
$ws.Range("I122").Value2 = 0
$ws.Range("K122").Value2 = 0
$ws.Range("H126").Value2 = 3108.2
$ws.Range("I126").Value2 = 3215.3333
$ws.Range("J126").Value2 = 3036.7778
$ws.Range("K126").Value2 = 9645.999899999999
$ws.Range("L126").Value2 = 9110.3334
$ws.Range("M126").Value2 = -7175.999899999999
$ws.Range("N126").Value2 = -14050.3334
$ws.Range("H132").Value2 = 4890.222
$ws.Range("I132").Value2 = 4666.6665
$ws.Range("J132").Value2 = 5002
$ws.Range("K132").Value2 = 13999.9995
$ws.Range("L132").Value2 = 15006
$ws.Range("M132").Value2 = -11469.9995
$ws.Range("N132").Value2 = -20066
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value2 = 30000
$ws.Range("J94").Value2 = 30000
$ws.Range("L94").Value2 = 30000
$ws.Range("N94").Value2 = -31352
$ws.Range("H132").Value2 = 4291.6665
$ws.Range("I132").Value2 = 4536.067
$ws.Range("J132").Value2 = 3069.6667
$ws.Range("K132").Value2 = 13608.201
$ws.Range("L132").Value2 = 9209.000100000001
$ws.Range("M132").Value2 = -11078.201
$ws.Range("N132").Value2 = -14269.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value2 = 167531.08
$ws.Range("I81").Value2 = 125903
$ws.Range("J81").Value2 = 250787.25
$ws.Range("K81").Value2 = 251806
$ws.Range("L81").Value2 = 501574.5
$ws.Range("M81").Value2 = -250745
$ws.Range("N81").Value2 = -503696.5
$ws.Range("H84").Value2 = 167531.08
$ws.Range("I84").Value2 = 125903
$ws.Range("J84").Value2 = 250787.25
$ws.Range("K84").Value2 = 1259030
$ws.Range("L84").Value2 = 2507872.5
$ws.Range("M84").Value2 = -1253726
$ws.Range("N84").Value2 = -2518480.5
$ws.Range("H119").Value2 = 39581.8
$ws.Range("J119").Value2 = 39581.8
$ws.Range("L119").Value2 = 39581.8
$ws.Range("N119").Value2 = -49257.8
